$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the existing 2019 column (I4:I18) formatting into the new 2020 column (J4:J18)
# so the new column inherits the same fonts / number formats / borders / alignment.
$ws.Range("I4:I18").Copy()
$ws.Range("J4").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# Header
$ws.Range("J4").Value = 2020

# Data rows
$ws.Range("J5").Value = 8017.9
$ws.Range("J6").Value = $null
$ws.Range("J7").Formula = "=J5-J8"
$ws.Range("J8").Value = 249.8
$ws.Range("J9").Value = $null
$ws.Range("J10").Value = 757.6
$ws.Range("J11").Value = 984.4
$ws.Range("J12").Value = 646.2
$ws.Range("J13").Value = 667.6
$ws.Range("J14").Value = 1147
$ws.Range("J15").Value = 961.1
$ws.Range("J16").Value = 2664.5
$ws.Range("J17").Value = 132.5
$ws.Range("J18").Value = 57

# Move the active-cell selection like the after-state of the workbook
$ws.Range("J19").Select() | Out-Null
